$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 258 (everything currently at/after row 258
# shifts down by one, e.g. old row 258 -> new row 259, old row 360 -> new row 361).
$ws.Rows("258:258").Insert()

# Populate the newly inserted row 258 with the new data record.
$ws.Range("A258").Value = 7
$ws.Range("B258").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C258").Value = "Ñuble"
$ws.Range("D258").Value = 45146
$ws.Range("E258").Value = 16
$ws.Range("F258").Value = "Fruta"
$ws.Range("G258").Value = 100101
$ws.Range("H258").Value = "Berries"
$ws.Range("I258").Value = 100101007
$ws.Range("J258").Value = "Kiwi"
$ws.Range("K258").Value = "Hayward"
$ws.Range("L258").Value = "Primera"
$ws.Range("M258").Value = 120
$ws.Range("N258").Value = 13000
$ws.Range("O258").Value = 14000
$ws.Range("P258").Value = 13500
$ws.Range("Q258").Value = "$/bandeja 18 kilos"
$ws.Range("R258").Value = "Región de O'Higgins"
$ws.Range("S258").Value = 750
$ws.Range("T258").Value = 18
